$d = $word.ActiveDocument

# Locate the exact run text that needs to be split: "Effective: October 31, 2025 | v1.0"
$r = $d.Content
$found = $r.Find.Execute("Effective: October 31, 2025 | v1.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Effective: October 31, 2025 | v1.0' text to edit."
}

# Remove the matched text; the Range collapses to the insertion point while
# staying inside the original paragraph (preserving its pPr / paraId / rsids).
$r.Delete()

# Re-insert the same visible text, but split into three runs so that the
# year "2025" is wrapped in proofErr gramStart/gramEnd markers (matching a
# Word grammar-check pass flagging the date as a potential grammar issue).
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7103BB4B" w14:textId="77777777" w:rsidR="0012355D" w:rsidRDefault="00000000"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Effective: October 31, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>2025</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> | v1.0</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
